$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (same set applied uniformly to every data row, B:I)
$values = @(0.353672031788087, -2.142122532649438, -0.5000481927095577, -0.8706523448567105, 0.7152945399284363, 0.6458885073661804, 0.7320140600204468, 0.6864180564880371)

# Data rows are 2 through 26; columns B (2) through I (9)
for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
